# React/Django - First login page
# Adds new TextContent rows (29-33 / row 30-34) for the "help topics" /
# "Save" / validation-message strings, widens column B to fit the longer
# text, and updates the active selection/scroll position to match the
# author's final view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the five new TextContent rows -------------------------------
$newRows = @(
    @(29, "On what issues can you help others?"),
    @(30, "What issues do you need help with?"),
    @(31, "Save"),
    @(32, "It is mandatory to fill this field."),
    @(33, "Must contain at least 50 and maximum 300 letters.")
)

$startRow = 30
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $id = $newRows[$i][0]
    $text = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 2).Value = $text
    $ws.Cells.Item($r, 3).Value = 1
}

# --- Column B needs to be a bit wider to comfortably fit the new,
#     longer strings (also drops the old "best fit" auto width) ----------
$ws.Columns.Item(2).ColumnWidth = 44.83

# --- Restore the view state left behind after editing: scrolled down
#     with B27 as the active cell ----------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("B27").Select()
